$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for rows 2-6: 45174 -> 45175 (2023-09-05 -> 2023-09-06)
$ws.Range("C2:C6").Value = 45175
